# Auto-generated Excel COM-interop edit script
# Applies numeric refresh to Leve profit-calculation sheets (H..N columns)
# per scheduled market-price runner update.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3725
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 3725
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 3725
$ws.Range("M64").Value = ""
$ws.Range("N64").Value = -4221

$ws.Range("H67").Value = 3725
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 3725
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 3725
$ws.Range("M67").Value = ""
$ws.Range("N67").Value = -5441

$ws.Range("H74").Value = 3680
$ws.Range("J74").Value = 3680
$ws.Range("L74").Value = 3680
$ws.Range("N74").Value = -5552

$ws.Range("H77").Value = 3680
$ws.Range("J77").Value = 3680
$ws.Range("L77").Value = 18400
$ws.Range("N77").Value = -27760

$ws.Range("H100").Value = 1919.8077
$ws.Range("I100").Value = 1117.3077
$ws.Range("J100").Value = 2722.3076
$ws.Range("K100").Value = 1117.3077
$ws.Range("L100").Value = 2722.3076
$ws.Range("M100").Value = -576.3077000000001
$ws.Range("N100").Value = -3804.3076

$ws.Range("H113").Value = 3721.6667
$ws.Range("I113").Value = 2693.0625
$ws.Range("J113").Value = 5217.8184
$ws.Range("K113").Value = 2693.0625
$ws.Range("L113").Value = 5217.8184
$ws.Range("M113").Value = 560.9375
$ws.Range("N113").Value = -11725.8184

$ws.Range("H127").Value = 71431880
$ws.Range("J127").Value = 4446.6
$ws.Range("L127").Value = 13339.8
$ws.Range("N127").Value = -23259.8

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = ""


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 56882.723
$ws.Range("I74").Value = 72677.57000000001
$ws.Range("J74").Value = 1600.75
$ws.Range("K74").Value = 72677.57000000001
$ws.Range("L74").Value = 1600.75
$ws.Range("M74").Value = -71803.57000000001
$ws.Range("N74").Value = -3348.75

$ws.Range("H77").Value = 56882.723
$ws.Range("I77").Value = 72677.57000000001
$ws.Range("J77").Value = 1600.75
$ws.Range("K77").Value = 363387.85
$ws.Range("L77").Value = 8003.75
$ws.Range("M77").Value = -359019.85
$ws.Range("N77").Value = -16739.75

$ws.Range("H110").Value = 1551.3334
$ws.Range("I110").Value = 1778.5385
$ws.Range("K110").Value = 1778.5385
$ws.Range("M110").Value = 266.4614999999999


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H47").Value = 98842
$ws.Range("J47").Value = 98842
$ws.Range("L47").Value = 98842
$ws.Range("N47").Value = -99882

$ws.Range("H99").Value = 1541.0344
$ws.Range("I99").Value = 1294.5
$ws.Range("J99").Value = 2088.889
$ws.Range("K99").Value = 1294.5
$ws.Range("L99").Value = 2088.889
$ws.Range("M99").Value = 203.5
$ws.Range("N99").Value = -5084.889

$ws.Range("H107").Value = 1729.4166
$ws.Range("I107").Value = 1624.6111
$ws.Range("J107").Value = 2043.8334
$ws.Range("K107").Value = 1624.6111
$ws.Range("L107").Value = 2043.8334
$ws.Range("M107").Value = 295.3888999999999
$ws.Range("N107").Value = -5883.8334


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 959.1429000000001
$ws.Range("I16").Value = 777.5
$ws.Range("J16").Value = 1031.8
$ws.Range("K16").Value = 777.5
$ws.Range("L16").Value = 1031.8
$ws.Range("M16").Value = -490.5
$ws.Range("N16").Value = -1605.8

$ws.Range("H28").Value = 25160.75
$ws.Range("J28").Value = 25160.75
$ws.Range("L28").Value = 25160.75
$ws.Range("N28").Value = -25650.75

$ws.Range("H31").Value = 22728834
$ws.Range("I31").Value = 30304088
$ws.Range("J31").Value = 3073
$ws.Range("K31").Value = 30304088
$ws.Range("L31").Value = 3073
$ws.Range("M31").Value = -30303793
$ws.Range("N31").Value = -3663

$ws.Range("H34").Value = 22728834
$ws.Range("I34").Value = 30304088
$ws.Range("J34").Value = 3073
$ws.Range("K34").Value = 30304088
$ws.Range("L34").Value = 3073
$ws.Range("M34").Value = -30303886
$ws.Range("N34").Value = -3477

$ws.Range("H92").Value = 29863.637
$ws.Range("J92").Value = 29863.637
$ws.Range("L92").Value = 29863.637
$ws.Range("N92").Value = -34855.637

$ws.Range("H107").Value = 1505.4546
$ws.Range("I107").Value = 514.6923
$ws.Range("J107").Value = 5185.4287
$ws.Range("K107").Value = 514.6923
$ws.Range("L107").Value = 5185.4287
$ws.Range("M107").Value = 1405.3077
$ws.Range("N107").Value = -9025.4287

$ws.Range("H113").Value = 959.1429000000001
$ws.Range("I113").Value = 777.5
$ws.Range("J113").Value = 1031.8
$ws.Range("K113").Value = 777.5
$ws.Range("L113").Value = 1031.8
$ws.Range("M113").Value = 1392.5
$ws.Range("N113").Value = -5371.8

$ws.Range("H122").Value = 1779
$ws.Range("I122").Value = 1742
$ws.Range("J122").Value = 1816
$ws.Range("K122").Value = 5226
$ws.Range("L122").Value = 5448
$ws.Range("M122").Value = -2776
$ws.Range("N122").Value = -10348

$ws.Range("H141").Value = 63467.715
$ws.Range("J141").Value = 63467.715
$ws.Range("L141").Value = 63467.715
$ws.Range("N141").Value = -73827.715


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 200
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").Value = ""

$ws.Range("H75").Value = 4223.091
$ws.Range("J75").Value = 4717.1113
$ws.Range("L75").Value = 14151.3339
$ws.Range("N75").Value = -16147.3339

$ws.Range("H78").Value = 4223.091
$ws.Range("J78").Value = 4717.1113
$ws.Range("L78").Value = 42454.00169999999
$ws.Range("N78").Value = -52438.00169999999

$ws.Range("H109").Value = 2975.2856
$ws.Range("J109").Value = 3360
$ws.Range("L109").Value = 10080
$ws.Range("N109").Value = -12160

$ws.Range("H115").Value = 2354.7144
$ws.Range("J115").Value = 3350
$ws.Range("L115").Value = 10050
$ws.Range("N115").Value = -12400

$ws.Range("H131").Value = 909.57574
$ws.Range("J131").Value = 909.57574
$ws.Range("L131").Value = 2728.72722
$ws.Range("N131").Value = -12808.72722


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5191.7144
$ws.Range("I70").Value = 5115.7896
$ws.Range("J70").Value = 5239.8
$ws.Range("K70").Value = 5115.7896
$ws.Range("L70").Value = 5239.8
$ws.Range("M70").Value = -4845.7896
$ws.Range("N70").Value = -5779.8

$ws.Range("H73").Value = 5191.7144
$ws.Range("I73").Value = 5115.7896
$ws.Range("J73").Value = 5239.8
$ws.Range("K73").Value = 5115.7896
$ws.Range("L73").Value = 5239.8
$ws.Range("M73").Value = -4179.7896
$ws.Range("N73").Value = -7111.8

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").Value = ""

$ws.Range("H122").Value = 31100.527
$ws.Range("I122").Value = 45946.086
$ws.Range("J122").Value = 4835.3076
$ws.Range("K122").Value = 137838.258
$ws.Range("L122").Value = 14505.9228
$ws.Range("M122").Value = -135388.258
$ws.Range("N122").Value = -19405.9228

$ws.Range("H126").Value = 2475.5
$ws.Range("I126").Value = 3140.5
$ws.Range("J126").Value = 1810.5
$ws.Range("K126").Value = 9421.5
$ws.Range("L126").Value = 5431.5
$ws.Range("M126").Value = -6951.5
$ws.Range("N126").Value = -10371.5


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1645.2
$ws.Range("I7").Value = 1496
$ws.Range("J7").Value = 1993.3334
$ws.Range("K7").Value = 1496
$ws.Range("L7").Value = 1993.3334
$ws.Range("M7").Value = -1384
$ws.Range("N7").Value = -2217.3334

$ws.Range("H68").Value = 7831.6313
$ws.Range("I68").Value = 18300.334
$ws.Range("J68").Value = 2999.923
$ws.Range("K68").Value = 18300.334
$ws.Range("L68").Value = 2999.923
$ws.Range("M68").Value = -17551.334
$ws.Range("N68").Value = -4497.923

$ws.Range("H71").Value = 7831.6313
$ws.Range("I71").Value = 18300.334
$ws.Range("J71").Value = 2999.923
$ws.Range("K71").Value = 91501.67
$ws.Range("L71").Value = 14999.615
$ws.Range("M71").Value = -87757.67
$ws.Range("N71").Value = -22487.615

$ws.Range("H108").Value = 29836
$ws.Range("J108").Value = 29836
$ws.Range("L108").Value = 29836
$ws.Range("N108").Value = -37516

$ws.Range("H122").Value = 12490
$ws.Range("I122").Value = 20000
$ws.Range("J122").Value = 4980
$ws.Range("K122").Value = 60000
$ws.Range("L122").Value = 14940
$ws.Range("M122").Value = -57550
$ws.Range("N122").Value = -19840

$ws.Range("H126").Value = 1645.2
$ws.Range("I126").Value = 1496
$ws.Range("J126").Value = 1993.3334
$ws.Range("K126").Value = 4488
$ws.Range("L126").Value = 5980.0002
$ws.Range("M126").Value = -2018
$ws.Range("N126").Value = -10920.0002


Write-Host "Applied Belias_Profits market-price refresh across ALC, ARM, BSM, CRP, CUL, GSM, LTW"